# Undo/Redo activity diagram: replace "address book" wording with
# "hotel management system" and adjust the affected shapes' geometry
# to match the re-flowed layout.
#
# EMU-exact positioning helper: PowerPoint COM's Shape.Left/Top/Width/
# Height are expressed in points (1 pt = 12700 EMU). Adding half an EMU
# before dividing compensates for truncation in the host's pt->EMU
# conversion so the stored <a:off>/<a:ext> values land on the exact EMU
# figures from the source file instead of being 1 EMU short.
function EMU([double]$emu) {
    return ($emu + 0.5) / 12700.0
}

function Get-ShapeById($shapes, $id) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.Id -eq $id) {
            return $sh
        }
    }
    return $null
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- "[command commits address book]" textbox (id 48) ---------------
$sh48 = Get-ShapeById $s.Shapes 48
$tr48 = $sh48.TextFrame.TextRange
# Run 1 is the literal "[", run 2 is "command commits address book]".
$run48_2 = $tr48.Characters(2, $tr48.Length - 1)
$run48_2.Text = "command commits  hotel management system]"

$sh48.Left   = EMU 4202312
$sh48.Top    = EMU 1910082
$sh48.Width  = EMU 2066045
$sh48.Height = EMU 923458

# --- "Purge redundant states ... address book ..." rounded rect (id 51) --
$sh51 = Get-ShapeById $s.Shapes 51
$tr51 = $sh51.TextFrame.TextRange
# Runs: "Purge redundant states and then save address book to "
#       "addressBookStateList"
#       " "
# Edit the later run first so the earlier run's character offsets stay valid.
$run51_2 = $tr51.Characters(54, 20)
$run51_2.Text = "hotelManagementSystemStateList"
$run51_1 = $tr51.Characters(1, 53)
$run51_1.Text = "Purge redundant states and then save hotel management system to "

$sh51.Left   = EMU 6237767
$sh51.Top    = EMU 2311019
$sh51.Width  = EMU 3568526
$sh51.Height = EMU 814659

# --- Elbow Connector 65 (id 66): re-routed after shape 51 moved ------
$sh66 = Get-ShapeById $s.Shapes 66
$sh66.Left   = EMU 5940526
$sh66.Top    = EMU 2775207
$sh66.Width  = EMU 354098
$sh66.Height = EMU 240383

# --- Elbow Connector 73 (id 74): re-routed after shape 51 moved ------
$sh74 = Get-ShapeById $s.Shapes 74
$sh74.Left   = EMU 9806293
$sh74.Top    = EMU 2718349
$sh74.Width  = EMU 84768
$sh74.Height = EMU 357048
